$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the SUM formula in O1 and replace with static value 11
$ws.Range("O1").Value = 11

# Add new values in column M (rows 2-11)
$ws.Range("M2").Value = 0
$ws.Range("M3").Value = 1
$ws.Range("M4").Value = "3 downto 2"
$ws.Range("M4").NumberFormatLocal = "General"
$ws.Range("M5").Value = 4
$ws.Range("M6").Value = 5
$ws.Range("M7").Value = 6
$ws.Range("M8").Value = 7
$ws.Range("M9").Value = 8
$ws.Range("M10").Value = 9
$ws.Range("M11").Value = 10

# Update selection/view
$ws.Range("O2").Select()
